# Sjoberg - RMed 2023 CRC Poster.pptx edits
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Update the poster title text (Shape Id 28, index 1 - "Title 1")
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Estimation of Patient Death Accounting for Competing Events"

# 2) Move the "Methods" header textbox up (Shape Id 35, index 7 - "Text Placeholder 5")
#    off x=1358582 y=13010188 -> y=12879560 (EMU); 914400 EMU per inch, 12700 EMU per point
$s.Shapes.Item(7).Top = 12879560 / 12700

# 3) Move the large "Methods" body textbox up to match (Shape Id 36, index 8 - "Text Placeholder 6")
#    off x=1356994 y=13863139 -> y=13732511 (EMU)
$s.Shapes.Item(8).Top = 13732511 / 12700

# 4) Shift "Picture 8" to the right (Shape Id 9, index 10)
#    off x=17136847 y=4552682 -> x=17561388 (EMU)
$s.Shapes.Item(10).Left = 17561388 / 12700
